# Append a new data row (row 44) to Sheet1, mirroring the existing rows:
#   A -> date text, B -> weekday text, C -> hour (number), D -> ranking (number)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds a date-like string ("2025/10/01") that must stay as literal
# text rather than being auto-converted to a date serial number by Excel's
# smart entry. Temporarily force a text number format, assign the value,
# then restore the cell to the default "Normal" style so no extra
# formatting is left behind on the cell (matching the other rows, which
# carry no style index).
$ws.Range("A44").NumberFormat = "@"
$ws.Range("A44").Value = "2025/10/01"
$ws.Range("A44").Style = "Normal"

$ws.Range("B44").Value = "水"
$ws.Range("C44").Value = 9
$ws.Range("D44").Value = 12
